$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) is stored as plain text in this workbook (values like
# "1.562.38" are not valid numbers, and even single-dot values such as "289.61" are
# meant to stay literal text). Force a text number format first so Excel does not
# reinterpret the new price strings as floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.343.65"
$ws.Range("E2").Value = "  -4.89%  "
$ws.Range("D3").Value = "1.562.06"
$ws.Range("E3").Value = "  -5.08%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "289.53"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("D7").Value = "0.3716"
$ws.Range("E7").Value = "  -2.07%  "
$ws.Range("D8").Value = "49.20"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").Value = "0.3400"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").Value = "1.165"
$ws.Range("E10").Value = "  -4.18%  "
$ws.Range("D11").Value = "0.07636"
$ws.Range("E11").Value = "  -5.33%  "
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "21.45"
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").Value = "6.050"
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").Value = "6.909"
$ws.Range("E15").Value = "  -4.64%  "
$ws.Range("D16").Value = "1.563.34"
$ws.Range("E16").Value = "  -4.95%  "
$ws.Range("E17").Value = "  -6.95%  "
$ws.Range("D18").Value = "90.02"
$ws.Range("E18").Value = "  -5.21%  "
$ws.Range("D19").Value = "0.06724"
$ws.Range("E19").Value = "  -3.55%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "6.229"
$ws.Range("E21").Value = "  -5.81%  "
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").Value = "0.5309"
$ws.Range("E23").Value = "  -7.25%  "
$ws.Range("E24").Value = "  -3.58%  "
$ws.Range("D25").Value = "22.349.03"
$ws.Range("E25").Value = "  -4.91%  "
$ws.Range("D26").Value = "2.399"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "2.829"
$ws.Range("E27").Value = "  -4.47%  "
$ws.Range("D28").Value = "20.16"
$ws.Range("E28").Value = "  -4.00%  "
$ws.Range("D29").Value = "145.52"
$ws.Range("E29").Value = "  -3.88%  "
$ws.Range("D30").Value = "4.980"
$ws.Range("E30").Value = "  -3.98%  "
$ws.Range("D31").Value = "125.35"
$ws.Range("E31").Value = "  -4.76%  "
$ws.Range("D32").Value = "1.742.00"
$ws.Range("E32").Value = "  -4.35%  "
$ws.Range("D33").Value = "6.198"
$ws.Range("E33").Value = "  -9.27%  "
$ws.Range("D34").Value = "2.009"
$ws.Range("E34").Value = "  -6.13%  "
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("D36").Value = "9.992"
$ws.Range("E36").Value = "  -10.55%  "
$ws.Range("D37").Value = "0.08437"
$ws.Range("E37").Value = "  -4.01%  "
$ws.Range("D38").Value = "0.02534"
$ws.Range("E38").Value = "  -5.74%  "
$ws.Range("D39").Value = "0.2320"
$ws.Range("E39").Value = "  -4.08%  "
$ws.Range("D40").Value = "5.525"
$ws.Range("E40").Value = "  -6.43%  "
$ws.Range("D41").Value = "0.06387"
$ws.Range("E41").Value = "  -5.76%  "
$ws.Range("D42").Value = "1.299"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "11.68"
$ws.Range("E43").Value = "  -8.76%  "
$ws.Range("D44").Value = "0.6335"
$ws.Range("E44").Value = "  -7.69%  "
$ws.Range("D45").Value = "14.07"
$ws.Range("E45").Value = "  -9.47%  "
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  -6.36%  "
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("E49").Value = "  -6.85%  "
$ws.Range("D50").Value = "1.269"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").Value = "124.47"
$ws.Range("E51").Value = "  -1.95%  "
